$wb = $excel.ActiveWorkbook

# --- Sheet "8_": convert A2:A6 from plain numbers to quote-prefixed text ---
# (values 2.5, 5, 7.5, 10, 12.5 stay the same, just re-typed as quoted text)
$ws8 = $wb.Worksheets.Item("8_")
$ws8.Range("A2").Value = "''2.5'"
$ws8.Range("A3").Value = "''5'"
$ws8.Range("A4").Value = "''7.5'"
$ws8.Range("A5").Value = "''10'"
$ws8.Range("A6").Value = "''12.5'"
$null = $ws8.Range("A2:A6").Select()

# --- Sheet "9_": same quoted-text values copied into A2:A6 ---
$ws9 = $wb.Worksheets.Item("9_")
$ws9.Range("A2").Value = "''2.5'"
$ws9.Range("A3").Value = "''5'"
$ws9.Range("A4").Value = "''7.5'"
$ws9.Range("A5").Value = "''10'"
$ws9.Range("A6").Value = "''12.5'"
$null = $ws9.Range("A2:A6").Select()

# --- Sheet "10_" becomes the active sheet/tab, with a new selection ---
$ws10 = $wb.Worksheets.Item("10_")
$null = $ws10.Range("B19").Select()
$null = $ws10.Activate()
